# Rebuild the stock-price table: 14 new ticker rows are inserted at the
# top of the data block (rows 2-15) and the 49 original rows shift down
# to rows 16-64. Column A keeps its running index (0-based) and its bold/
# bordered header-row style; columns B-H keep the default (unstyled) look.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object 'object[,]' 63,7
$data[0,0] = "BANKINDIA.NS"
$data[0,1] = [double]110.2200012207031
$data[0,2] = [double]104.1999969482422
$data[0,3] = [double]111.2200012207031
$data[0,4] = [double]104.1999969482422
$data[0,5] = [double]107.8600006103516
$data[0,6] = "BANKINDIA"
$data[1,0] = "COCHINSHIP.NS"
$data[1,1] = [double]1484.599975585938
$data[1,2] = [double]1426.199951171875
$data[1,3] = [double]1492
$data[1,4] = [double]1417.900024414062
$data[1,5] = [double]1443.400024414062
$data[1,6] = "COCHINSHIP"
$data[2,0] = "EXCELINDUS.NS"
$data[2,1] = [double]1024
$data[2,2] = [double]1004.5
$data[2,3] = [double]1038
$data[2,4] = [double]991.0999755859375
$data[2,5] = [double]1046.400024414062
$data[2,6] = "EXCELINDUS"
$data[3,0] = "EXIDEIND.NS"
$data[3,1] = [double]360.9500122070312
$data[3,2] = [double]350.5
$data[3,3] = [double]362.2000122070312
$data[3,4] = [double]350.5
$data[3,5] = [double]360.8500061035156
$data[3,6] = "EXIDEIND"
$data[4,0] = "HAPPSTMNDS.NS"
$data[4,1] = [double]577.5
$data[4,2] = [double]563.2999877929688
$data[4,3] = [double]580
$data[4,4] = [double]560.5499877929688
$data[4,5] = [double]585.0499877929688
$data[4,6] = "HAPPSTMNDS"
$data[5,0] = "HINDUNILVR.BO"
$data[5,1] = [double]2333.949951171875
$data[5,2] = [double]2324.10009765625
$data[5,3] = [double]2340
$data[5,4] = [double]2302
$data[5,5] = [double]2355.25
$data[5,6] = "HINDUNILVR"
$data[6,0] = "IDFCFIRSTB.BO"
$data[6,1] = [double]66.20999908447266
$data[6,2] = [double]64.84999847412109
$data[6,3] = [double]66.90000152587891
$data[6,4] = [double]64.84999847412109
$data[6,5] = [double]66.44999694824219
$data[6,6] = "IDFCFIRSTB"
$data[7,0] = "KPIGREEN.BO"
$data[7,1] = [double]348.6000061035156
$data[7,2] = [double]335.6000061035156
$data[7,3] = [double]351.75
$data[7,4] = [double]335.5499877929688
$data[7,5] = [double]358.8999938964844
$data[7,6] = "KPIGREEN"
$data[8,0] = "MAHABANK.BO"
$data[8,1] = [double]48.34999847412109
$data[8,2] = [double]47.13999938964844
$data[8,3] = [double]49.13000106811523
$data[8,4] = [double]47.13999938964844
$data[8,5] = [double]48.91999816894531
$data[8,6] = "MAHABANK"
$data[9,0] = "RELIANCE.BO"
$data[9,1] = [double]1377.75
$data[9,2] = [double]1372.050048828125
$data[9,3] = [double]1395
$data[9,4] = [double]1372.050048828125
$data[9,5] = [double]1404.849975585938
$data[9,6] = "RELIANCE"
$data[10,0] = "SAIL.NS"
$data[10,1] = [double]109.0100021362305
$data[10,2] = [double]107
$data[10,3] = [double]109.3199996948242
$data[10,4] = [double]106.25
$data[10,5] = [double]109.2399978637695
$data[10,6] = "SAIL"
$data[11,0] = "SOUTHBANK.BO"
$data[11,1] = [double]24.34000015258789
$data[11,2] = [double]23.79000091552734
$data[11,3] = [double]24.43000030517578
$data[11,4] = [double]23.79000091552734
$data[11,5] = [double]24.20000076293945
$data[11,6] = "SOUTHBANK"
$data[12,0] = "TATAMOTORS.BO"
$data[12,1] = [double]708.5
$data[12,2] = [double]662.5999755859375
$data[12,3] = [double]710.8499755859375
$data[12,4] = [double]662.5999755859375
$data[12,5] = [double]681.9000244140625
$data[12,6] = "TATAMOTORS"
$data[13,0] = "TCS.BO"
$data[13,1] = [double]3442.199951171875
$data[13,2] = [double]3380
$data[13,3] = [double]3445.14990234375
$data[13,4] = [double]3380
$data[13,5] = [double]3447.449951171875
$data[13,6] = "TCS"
$data[14,0] = "BAJAJHIND.BO"
$data[14,1] = [double]19.06999969482422
$data[14,2] = [double]17.96999931335449
$data[14,3] = [double]19.35000038146973
$data[14,4] = [double]17.96999931335449
$data[14,5] = [double]18.73999977111816
$data[14,6] = "BAJAJHIND"
$data[15,0] = "BANDHANBNK.BO"
$data[15,1] = [double]156.9499969482422
$data[15,2] = [double]153.6000061035156
$data[15,3] = [double]157.9499969482422
$data[15,4] = [double]153.0500030517578
$data[15,5] = [double]157.3000030517578
$data[15,6] = "BANDHANBNK"
$data[16,0] = "BANKBARODA.NS"
$data[16,1] = [double]220.0899963378906
$data[16,2] = [double]212.9900054931641
$data[16,3] = [double]220.5
$data[16,4] = [double]212.5500030517578
$data[16,5] = [double]217.2700042724609
$data[16,6] = "BANKBARODA"
$data[17,0] = "BHEL.BO"
$data[17,1] = [double]216.75
$data[17,2] = [double]210.3999938964844
$data[17,3] = [double]217.9499969482422
$data[17,4] = [double]210.3999938964844
$data[17,5] = [double]217.3500061035156
$data[17,6] = "BHEL"
$data[18,0] = "BIGBLOC.NS"
$data[18,1] = [double]62.11000061035156
$data[18,2] = [double]60.16999816894531
$data[18,3] = [double]62.84999847412109
$data[18,4] = [double]59.95999908447266
$data[18,5] = [double]62.84999847412109
$data[18,6] = "BIGBLOC"
$data[19,0] = "CANBK.BO"
$data[19,1] = [double]97.65000152587891
$data[19,2] = [double]93.66000366210938
$data[19,3] = [double]98.40000152587891
$data[19,4] = [double]93.66000366210938
$data[19,5] = [double]95.37999725341797
$data[19,6] = "CANBK"
$data[20,0] = "COCHINSHIP.BO"
$data[20,1] = [double]1484.400024414062
$data[20,2] = [double]1435
$data[20,3] = [double]1491.699951171875
$data[20,4] = [double]1415.449951171875
$data[20,5] = [double]1444.099975585938
$data[20,6] = "COCHINSHIP"
$data[21,0] = "DELTACORP.BO"
$data[21,1] = [double]84.80000305175781
$data[21,2] = [double]81.61000061035156
$data[21,3] = [double]87.75
$data[21,4] = [double]81.61000061035156
$data[21,5] = [double]88.54000091552734
$data[21,6] = "DELTACORP"
$data[22,0] = "EMBDL.NS"
$data[22,1] = [double]91.29000091552734
$data[22,2] = [double]90
$data[22,3] = [double]92.08999633789062
$data[22,4] = [double]89.22000122070312
$data[22,5] = [double]91.55999755859375
$data[22,6] = "EMBDL"
$data[23,0] = "EXIDEIND.NS"
$data[23,1] = [double]360.9500122070312
$data[23,2] = [double]350.5
$data[23,3] = [double]362.2000122070312
$data[23,4] = [double]350.5
$data[23,5] = [double]360.8500061035156
$data[23,6] = "EXIDEIND"
$data[24,0] = "GLENMARK.NS"
$data[24,1] = [double]1398.199951171875
$data[24,2] = [double]1395
$data[24,3] = [double]1414.199951171875
$data[24,4] = [double]1382.900024414062
$data[24,5] = [double]1418.900024414062
$data[24,6] = "GLENMARK"
$data[25,0] = "HDFCBANK.NS"
$data[25,1] = [double]1889.699951171875
$data[25,2] = [double]1910.099975585938
$data[25,3] = [double]1919
$data[25,4] = [double]1886.800048828125
$data[25,5] = [double]1928.5
$data[25,6] = "HDFCBANK"
$data[26,0] = "ICICIBANK.NS"
$data[26,1] = [double]1388.900024414062
$data[26,2] = [double]1415.199951171875
$data[26,3] = [double]1419.900024414062
$data[26,4] = [double]1387
$data[26,5] = [double]1435.5
$data[26,6] = "ICICIBANK"
$data[27,0] = "IDBI.NS"
$data[27,1] = [double]76.30000305175781
$data[27,2] = [double]75.5
$data[27,3] = [double]76.66999816894531
$data[27,4] = [double]75
$data[27,5] = [double]76.87000274658203
$data[27,6] = "IDBI"
$data[28,0] = "IDEA.BO"
$data[28,1] = [double]6.710000038146973
$data[28,2] = [double]6.519999980926514
$data[28,3] = [double]6.75
$data[28,4] = [double]6.460000038146973
$data[28,5] = [double]6.690000057220459
$data[28,6] = "IDEA"
$data[29,0] = "IDFCFIRSTB.BO"
$data[29,1] = [double]66.20999908447266
$data[29,2] = [double]64.84999847412109
$data[29,3] = [double]66.90000152587891
$data[29,4] = [double]64.84999847412109
$data[29,5] = [double]66.44999694824219
$data[29,6] = "IDFCFIRSTB"
$data[30,0] = "IEX.NS"
$data[30,1] = [double]189.3800048828125
$data[30,2] = [double]182
$data[30,3] = [double]190.0500030517578
$data[30,4] = [double]182
$data[30,5] = [double]189.8999938964844
$data[30,6] = "IEX"
$data[31,0] = "ITC.NS"
$data[31,1] = [double]423.5499877929688
$data[31,2] = [double]425.9500122070312
$data[31,3] = [double]429.1499938964844
$data[31,4] = [double]423
$data[31,5] = [double]430.6000061035156
$data[31,6] = "ITC"
$data[32,0] = "JPPOWER.BO"
$data[32,1] = [double]13.27999973297119
$data[32,2] = [double]12.60999965667725
$data[32,3] = [double]13.39999961853027
$data[32,4] = [double]12.60999965667725
$data[32,5] = [double]13.26000022888184
$data[32,6] = "JPPOWER"
$data[33,0] = "KALAMANDIR.NS"
$data[33,1] = [double]113.5100021362305
$data[33,2] = [double]113
$data[33,3] = [double]114.8899993896484
$data[33,4] = [double]113
$data[33,5] = [double]114.7799987792969
$data[33,6] = "KALAMANDIR"
$data[34,0] = "KPIGREEN.NS"
$data[34,1] = [double]348.3500061035156
$data[34,2] = [double]336
$data[34,3] = [double]352
$data[34,4] = [double]336
$data[34,5] = [double]359.25
$data[34,6] = "KPIGREEN"
$data[35,0] = "LTF.NS"
$data[35,1] = [double]161.5299987792969
$data[35,2] = [double]159.7400054931641
$data[35,3] = [double]164.1699981689453
$data[35,4] = [double]159.6600036621094
$data[35,5] = [double]163.5299987792969
$data[35,6] = "LTF"
$data[36,0] = "LAURUSLABS.NS"
$data[36,1] = [double]588.7999877929688
$data[36,2] = [double]578
$data[36,3] = [double]590
$data[36,4] = [double]572.25
$data[36,5] = [double]585.9500122070312
$data[36,6] = "LAURUSLABS"
$data[37,0] = "NTPC.NS"
$data[37,1] = [double]334.75
$data[37,2] = [double]334.2000122070312
$data[37,3] = [double]338.75
$data[37,4] = [double]330.6000061035156
$data[37,5] = [double]340
$data[37,6] = "NTPC"
$data[38,0] = "ONGC.NS"
$data[38,1] = [double]234.9600067138672
$data[38,2] = [double]229.1000061035156
$data[38,3] = [double]235.3000030517578
$data[38,4] = [double]228.4499969482422
$data[38,5] = [double]233.2200012207031
$data[38,6] = "ONGC"
$data[39,0] = "PNB.NS"
$data[39,1] = [double]91.97000122070312
$data[39,2] = [double]90
$data[39,3] = [double]92.34999847412109
$data[39,4] = [double]89.44999694824219
$data[39,5] = [double]91.36000061035156
$data[39,6] = "PNB"
$data[40,0] = "PREMIERENE.NS"
$data[40,1] = [double]945.6500244140625
$data[40,2] = [double]922
$data[40,3] = [double]954
$data[40,4] = [double]917.0499877929688
$data[40,5] = [double]944.7000122070312
$data[40,6] = "PREMIERENE"
$data[41,0] = "RPOWER.BO"
$data[41,1] = [double]38.65000152587891
$data[41,2] = [double]37.13999938964844
$data[41,3] = [double]38.84999847412109
$data[41,4] = [double]37
$data[41,5] = [double]38.29000091552734
$data[41,6] = "RPOWER"
$data[42,0] = "SBIN.NS"
$data[42,1] = [double]779.25
$data[42,2] = [double]756.5
$data[42,3] = [double]781.7000122070312
$data[42,4] = [double]755.5
$data[42,5] = [double]769
$data[42,6] = "SBIN"
$data[43,0] = "SUVEN.NS"
$data[43,1] = [double]142.4499969482422
$data[43,2] = [double]136.9900054931641
$data[43,3] = [double]144.0299987792969
$data[43,4] = [double]131.9900054931641
$data[43,5] = [double]137.8300018310547
$data[43,6] = "SUVEN"
$data[44,0] = "TATAMOTORS.NS"
$data[44,1] = [double]708.5
$data[44,2] = [double]667.5
$data[44,3] = [double]711
$data[44,4] = [double]666
$data[44,5] = [double]682.8499755859375
$data[44,6] = "TATAMOTORS"
$data[45,0] = "TECHM.NS"
$data[45,1] = [double]1493.699951171875
$data[45,2] = [double]1480
$data[45,3] = [double]1501.300048828125
$data[45,4] = [double]1478.099975585938
$data[45,5] = [double]1502.800048828125
$data[45,6] = "TECHM"
$data[46,0] = "TFCILTD.BO"
$data[46,1] = [double]196.3000030517578
$data[46,2] = [double]195
$data[46,3] = [double]198.8500061035156
$data[46,4] = [double]185
$data[46,5] = [double]198.1999969482422
$data[46,6] = "TFCILTD"
$data[47,0] = "UCOBANK.BO"
$data[47,1] = [double]29.59000015258789
$data[47,2] = [double]29.36000061035156
$data[47,3] = [double]29.89999961853027
$data[47,4] = [double]29.04999923706055
$data[47,5] = [double]29.98999786376953
$data[47,6] = "UCOBANK"
$data[48,0] = "VAKRANGEE.NS"
$data[48,1] = [double]8.840000152587891
$data[48,2] = [double]9
$data[48,3] = [double]9.180000305175781
$data[48,4] = [double]8.819999694824219
$data[48,5] = [double]9.229999542236328
$data[48,6] = "VAKRANGEE"
$data[49,0] = "VISAKAIND.NS"
$data[49,1] = [double]58.20000076293945
$data[49,2] = [double]57.0099983215332
$data[49,3] = [double]58.75
$data[49,4] = [double]57.0099983215332
$data[49,5] = [double]59.47999954223633
$data[49,6] = "VISAKAIND"
$data[50,0] = "VMM.NS"
$data[50,1] = [double]120.7900009155273
$data[50,2] = [double]118.5999984741211
$data[50,3] = [double]121.870002746582
$data[50,4] = [double]116.8499984741211
$data[50,5] = [double]119.1900024414062
$data[50,6] = "VMM"
$data[51,0] = "VOLTAS.BO"
$data[51,1] = [double]1235.75
$data[51,2] = [double]1185.800048828125
$data[51,3] = [double]1239.400024414062
$data[51,4] = [double]1185.800048828125
$data[51,5] = [double]1224.300048828125
$data[51,6] = "VOLTAS"
$data[52,0] = "WABAG.NS"
$data[52,1] = [double]1273.199951171875
$data[52,2] = [double]1250
$data[52,3] = [double]1282
$data[52,4] = [double]1250
$data[52,5] = [double]1292.900024414062
$data[52,6] = "WABAG"
$data[53,0] = "WIPRO.NS"
$data[53,1] = [double]242.0099945068359
$data[53,2] = [double]236.6000061035156
$data[53,3] = [double]242.9499969482422
$data[53,4] = [double]236.5
$data[53,5] = [double]241.5700073242188
$data[53,6] = "WIPRO"
$data[54,0] = "YESBANK.NS"
$data[54,1] = [double]20.02000045776367
$data[54,2] = [double]17.79999923706055
$data[54,3] = [double]20.3700008392334
$data[54,4] = [double]17.77000045776367
$data[54,5] = [double]18.22999954223633
$data[54,6] = "YESBANK"
$data[55,0] = "AWHCL.NS"
$data[55,1] = [double]503.7000122070312
$data[55,2] = [double]497
$data[55,3] = [double]512.9500122070312
$data[55,4] = [double]495.6499938964844
$data[55,5] = [double]509.3500061035156
$data[55,6] = "AWHCL"
$data[56,0] = "BANKINDIA.NS"
$data[56,1] = [double]110.2200012207031
$data[56,2] = [double]104.1999969482422
$data[56,3] = [double]111.2200012207031
$data[56,4] = [double]104.1999969482422
$data[56,5] = [double]107.8600006103516
$data[56,6] = "BANKINDIA"
$data[57,0] = "BSOFT.NS"
$data[57,1] = [double]384.4500122070312
$data[57,2] = [double]376.5
$data[57,3] = [double]386
$data[57,4] = [double]376.2999877929688
$data[57,5] = [double]385.9500122070312
$data[57,6] = "BSOFT"
$data[58,0] = "CONCOR.NS"
$data[58,1] = [double]648.6500244140625
$data[58,2] = [double]645.2999877929688
$data[58,3] = [double]659.1500244140625
$data[58,4] = [double]645.2999877929688
$data[58,5] = [double]657.4000244140625
$data[58,6] = "CONCOR"
$data[59,0] = "GAIL.NS"
$data[59,1] = [double]181.6000061035156
$data[59,2] = [double]178.7599945068359
$data[59,3] = [double]182.8999938964844
$data[59,4] = [double]178.6000061035156
$data[59,5] = [double]184.2700042724609
$data[59,6] = "GAIL"
$data[60,0] = "HINDUNILVR.BO"
$data[60,1] = [double]2333.949951171875
$data[60,2] = [double]2324.10009765625
$data[60,3] = [double]2340
$data[60,4] = [double]2302
$data[60,5] = [double]2355.25
$data[60,6] = "HINDUNILVR"
$data[61,0] = "RELIANCE.NS"
$data[61,1] = [double]1377.199951171875
$data[61,2] = [double]1385.5
$data[61,3] = [double]1394.800048828125
$data[61,4] = [double]1374.5
$data[61,5] = [double]1407
$data[61,6] = "RELIANCE"
$data[62,0] = "SOUTHBANK.BO"
$data[62,1] = [double]24.34000015258789
$data[62,2] = [double]23.79000091552734
$data[62,3] = [double]24.43000030517578
$data[62,4] = [double]23.79000091552734
$data[62,5] = [double]24.20000076293945
$data[62,6] = "SOUTHBANK"

$ws.Range("B2:H64").Value = $data

# Column A: sequential 0-based counter for every data row
$idx = New-Object 'object[,]' 63,1
for ($i = 0; $i -lt 63; $i++) { $idx[$i,0] = [double]$i }
$ws.Range("A2:A64").Value = $idx

# Extend column A's bold/bordered style (already present on A2:A50) down
# to the newly-added rows 51:64 without minting a brand-new cell style.
$ws.Range("A2").Copy()
$ws.Range("A51:A64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

